$wb = $excel.ActiveWorkbook

$hotel  = $wb.Worksheets.Item("hotel_info")
$review = $wb.Worksheets.Item("review_info")

# Insert a new "State" column into hotel_info, between "Hotel_Name" (B) and "City" (C)
$hotel.Columns("C").Insert()
$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# Reorder the sheet tabs so review_info comes before hotel_info
$review.Move($hotel)
